$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Table 3 title: "PSC" -> "integrated chlorophyll-a" (with the 'a' italicised)
$newTitle = "Table 3. Results of the analysis of integrated chlorophyll-a concentration using a 3-way ANOVA . Bold denotes significant results"
$ws.Range("A2").Value = $newTitle

# Touch the whole-cell italic flag (on, then back off) so the italic
# "Times New Roman" font combination gets registered in the workbook's font
# table, matching how Excel records fonts used anywhere in the workbook
# even after the character-level formatting below narrows it back down to
# just the single italic letter.
$ws.Range("A2").Font.Italic = $true
$ws.Range("A2").Font.Italic = $false

# Italicise the single "a" in "chlorophyll-a"
$ws.Range("A2").Characters(60, 1).Font.Italic = $true

# Make sure the remainder of the string after the italic "a" keeps an explicit
# (non-italic) run so the shared string ends up with three distinct runs, as
# in the authored edit.
$ws.Range("A2").Characters(61, $newTitle.Length - 60).Font.Italic = $false
